# Update cryptos list values per commit diff (Tue Nov 26 05:51:28 UTC 2024, GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay plain TEXT even when it looks numeric
# (e.g. "42.18"), matching the source data which is text, not a real number.
# Toggling NumberFormat to Text ("@") before the write, then resetting the
# cell Style back to "Normal" afterwards, keeps the cell free of any leftover
# number-format styling while the stored value remains a text string.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '94.611.32'
$ws.Range("E2").Value = '  -3.56%  '
$ws.Range("D3").Value = '3.428.74'
$ws.Range("E3").Value = '  +1.15%  '
$ws.Range("E4").Value = '  -0.01%  '
Set-TextValue "D5" '238.28'
$ws.Range("E5").Value = '  -6.17%  '
Set-TextValue "D6" '643.60'
$ws.Range("E6").Value = '  -2.78%  '
$ws.Range("E7").Value = '  -1.11%  '
$ws.Range("E8").Value = '  -4.33%  '
$ws.Range("E9").Value = '  +0.09%  '
Set-TextValue "D10" '0.985'
$ws.Range("E10").Value = '  -6.40%  '
$ws.Range("D11").Value = '3.424.86'
$ws.Range("E11").Value = '  +1.12%  '
$ws.Range("E12").Value = '  -4.24%  '
Set-TextValue "D13" '42.18'
$ws.Range("E13").Value = '  +0.28%  '
$ws.Range("E14").Value = '  +1.96%  '
$ws.Range("D15").Value = '94.360.26'
$ws.Range("E15").Value = '  -3.49%  '
$ws.Range("D16").Value = '4.067.86'
$ws.Range("E16").Value = '  +1.30%  '
$ws.Range("E17").Value = '  -1.11%  '
Set-TextValue "D18" '8.42'
$ws.Range("E18").Value = '  -7.25%  '
$ws.Range("D19").Value = '3.435.00'
Set-TextValue "D20" '17.59'
$ws.Range("E20").Value = '  -3.33%  '
Set-TextValue "D21" '11.69'
$ws.Range("E21").Value = '  +6.24%  '
$ws.Range("E22").Value = '  -5.27%  '
Set-TextValue "D23" '501.65'
$ws.Range("E23").Value = '  -2.11%  '
Set-TextValue "D24" '3.23'
$ws.Range("E24").Value = '  -6.37%  '
$ws.Range("E25").Value = '  -3.87%  '
Set-TextValue "D26" '6.57'
$ws.Range("E26").Value = '  -5.78%  '
Set-TextValue "D27" '94.49'
$ws.Range("E27").Value = '  -2.52%  '
$ws.Range("E28").Value = '  -3.11%  '
$ws.Range("D29").Value = '3.611.78'
$ws.Range("E29").Value = '  +1.20%  '
Set-TextValue "D30" '11.81'
$ws.Range("E30").Value = '  +3.01%  '
Set-TextValue "D31" '0.999'
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("E32").Value = '  +6.27%  '
$ws.Range("E33").Value = '  -3.11%  '
$ws.Range("B34").Value = 'Cronos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D34" '0.180'
$ws.Range("E34").Value = '  -3.98%  '
$ws.Range("B35").Value = 'Binance-PegBSC-USD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue "D35" '0.999'
$ws.Range("E35").Value = '  -0.09%  '
Set-TextValue "D36" '29.77'
$ws.Range("E36").Value = '  +2.69%  '
Set-TextValue "D37" '0.552'
$ws.Range("E37").Value = '  -2.14%  '
Set-TextValue "D38" '567.85'
$ws.Range("E38").Value = '  +6.20%  '
Set-TextValue "D39" '7.70'
$ws.Range("E39").Value = '  -3.64%  '
$ws.Range("E40").Value = '  -2.49%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("E42").Value = '  -1.29%  '
Set-TextValue "D43" '0.907'
$ws.Range("E43").Value = '  +5.77%  '
$ws.Range("E45").Value = '  -0.59%  '
Set-TextValue "D46" '3.72'
$ws.Range("E46").Value = '  +0.84%  '
Set-TextValue "D47" '5.70'
$ws.Range("E47").Value = '  +1.24%  '
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D48" '3.36'
$ws.Range("E48").Value = '  +4.19%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D49" '0.0413'
$ws.Range("E49").Value = '  -4.40%  '
Set-TextValue "D50" '54.90'
$ws.Range("E50").Value = '  -2.05%  '
$ws.Range("E51").Value = '  -2.57%  '
